$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D-column (Price) cells to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00", "23.99") into numbers,
# then restore the default "Normal" style so no extra formatting is applied.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) values per row
$ws.Range("D2").Value = '59.417.05'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.345.69'
$ws.Range("E3").Value = '  -3.01%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '559.25'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = '131.82'
$ws.Range("E6").Value = '  -3.94%  '
$ws.Range("E8").Value = '  -3.11%  '
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").Value = '0.339'
$ws.Range("E12").Value = '  -4.02%  '
$ws.Range("D13").Value = '23.99'
$ws.Range("E13").Value = '  -4.85%  '
$ws.Range("D14").Value = '2.765.67'
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").Value = '59.408.76'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D17").Value = '2.339.28'
$ws.Range("E17").Value = '  -3.42%  '
$ws.Range("D19").Value = '4.42'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '318.31'
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("D21").Value = '6.59'
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '63.75'
$ws.Range("E23").Value = '  -3.24%  '
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").Value = '8.35'
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("E27").Value = '  -3.03%  '
$ws.Range("D29").Value = '170.97'
$ws.Range("E29").Value = '  +1.63%  '
$ws.Range("D30").Value = '0.0₃0743'
$ws.Range("E30").Value = '  -4.39%  '
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("E32").Value = '  +4.67%  '
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("D34").Value = '17.97'
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("E40").Value = '  -2.59%  '
$ws.Range("D41").Value = '311.19'
$ws.Range("E41").Value = '  -4.39%  '
$ws.Range("D42").Value = '144.73'
$ws.Range("E42").Value = '  +3.10%  '
$ws.Range("D43").Value = '3.48'
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("E44").Value = '  -1.77%  '
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("E46").Value = '  -2.46%  '
$ws.Range("D47").Value = '18.71'
$ws.Range("E47").Value = '  -4.64%  '
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("D49").Value = '11.07'
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("E51").Value = '  -0.42%  '

# Restore default style on the D cells we touched (removes the temporary
# text number-format so the saved file keeps using the original default style)
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
